# Update the marksheet's "Total" row: corrected/total marks.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Marking row (B11): number right used in score calc
$ws.Range("B11").Value = 5

# Total row (B12): total correct marks, and E12 text "corr/total"
$ws.Range("B12").Value = 95
$ws.Range("E12").Value = "95/140"
